$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new columns R, S, T, U for the 11-month environment --------------
# R: copy of column A (SQRT(B*C))
# S: new 11-month fluence formula with polarization correction factor
# T: copy of the O-column style formula (C-B)
# U: copy of the P-column style formula (=B)

# Row 2 gets standalone formulas (matches the source workbook's row-2 pattern)
$ws.Range("R2").Formula = "=A2"
$ws.Range("S2").Formula = "=E2*(11/12)*365*24*3600*(1/(0.984)^2)"
$ws.Range("T2").Formula = "=C2-B2"
$ws.Range("U2").Formula = "=B2"

# Rows 3:29 get filled as one shared-formula block each (matches A3:A29 etc.)
$ws.Range("R3:R29").Formula = "=A3"
$ws.Range("S3:S29").Formula = "=E3*(11/12)*365*24*3600*(1/(0.984)^2)"
$ws.Range("T3:T29").Formula = "=C3-B3"
$ws.Range("U3:U29").Formula = "=B3"

# Apply the same number format/style used by the other computed columns (s="2")
$ws.Range("R2:U29").NumberFormat = $ws.Range("M2").NumberFormat

# --- View state -------------------------------------------------------------
$ws.Range("S3").Select()
$excel.ActiveWindow.ScrollColumn = 3
